# Re-generate the "haul" time strings: pad single-digit minutes/seconds
# with a leading zero, e.g. "198 ч. 23 мин. 8 сек." -> "198 ч. 23 мин. 08 сек."
# Hours are left untouched (they are not zero-padded in the source data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow  = $firstRow + $used.Rows.Count - 1
$firstCol = $used.Column
$lastCol  = $firstCol + $used.Columns.Count - 1

# Find the "Общее время" (total time) column from the header row.
$timeCol = 0
for ($c = $firstCol; $c -le $lastCol; $c++) {
    $header = $ws.Cells.Item($firstRow, $c).Value()
    if ($header -eq "Общее время") {
        $timeCol = $c
        break
    }
}
if ($timeCol -eq 0) {
    # Fallback: column D, matching the known report layout.
    $timeCol = 4
}

$pattern = '^(\d+) ч\. (\d+) мин\. (\d+) сек\.$'

for ($r = $firstRow + 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $timeCol)
    $text = $cell.Value()
    if ($null -eq $text) {
        continue
    }
    if ($text -match $pattern) {
        $hours   = $matches[1]
        $minutes = $matches[2]
        $seconds = $matches[3]
        if ($minutes.Length -lt 2 -or $seconds.Length -lt 2) {
            $minutes2 = $minutes.PadLeft(2, '0')
            $seconds2 = $seconds.PadLeft(2, '0')
            $cell.Value = "$hours ч. $minutes2 мин. $seconds2 сек."
        }
    }
}
